$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a Text number format on cells whose new values look like numbers,
# so Excel keeps them as text (matching the sheet's existing text-based
# price/volume columns) instead of silently converting them to numerics.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Assign the updated cell values
$ws.Range("D2").Value = "63.073.14"
$ws.Range("D3").Value = "2.953.31"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "595.11"
$ws.Range("D6").Value = "149.27"
$ws.Range("E6").Value = "  +3.12%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "2.949.45"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("D10").Value = "7.14"
$ws.Range("E10").Value = "  +3.09%  "
$ws.Range("D11").Value = "0.150"
$ws.Range("E11").Value = "  +6.49%  "
$ws.Range("D12").Value = "0.441"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").Value = "0.0000235"
$ws.Range("E13").Value = "  +5.16%  "
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").Value = "3.444.51"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").Value = "63.041.57"
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("D18").Value = "6.70"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "2.957.78"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("D20").Value = "442.15"
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("D21").Value = "13.50"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "0.668"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").Value = "7.02"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("D24").Value = "11.21"
$ws.Range("E24").Value = "  +2.88%  "
$ws.Range("D25").Value = "80.98"
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("D26").Value = "2.13"
$ws.Range("E26").Value = "  -1.86%  "
$ws.Range("D27").Value = "11.78"
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").Value = "7.33"
$ws.Range("E29").Value = "  +6.39%  "
$ws.Range("D30").Value = "2.23"
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("E32").Value = "  +16.07%  "
$ws.Range("D33").Value = "26.47"
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("E34").Value = "  -0.71%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "0.990"
$ws.Range("E36").Value = "  -1.65%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").Value = "5.61"
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "3.09"
$ws.Range("E38").Value = "  +3.72%  "
$ws.Range("E39").Value = "  +2.38%  "
$ws.Range("D40").Value = "49.78"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").Value = "8.51"
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").Value = "  -4.24%  "
$ws.Range("D43").Value = "0.280"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").Value = "39.14"
$ws.Range("E44").Value = "  -7.93%  "
$ws.Range("D45").Value = "2.702.02"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").Value = "135.62"
$ws.Range("E46").Value = "  +1.44%  "
$ws.Range("D47").Value = "0.0338"
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("D48").Value = "361.46"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "0.104"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").Value = "22.80"
$ws.Range("E51").Value = "  -2.95%  "
